$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.083188353011337
$ws.Range("C2").Value = 0.157306099068677
$ws.Range("D2").Value = 0.05162501809898856
$ws.Range("F2").Value = 6.672031919112214
$ws.Range("G2").Value = 0.002680968051116933
$ws.Range("J2").Value = 0.1774384546004057
$ws.Range("L2").Value = 0.2849453852701416
$ws.Range("M2").Value = 0.5816497773916325

$ws.Range("B3").Value = 3.043885440472025
$ws.Range("C3").Value = 0.1439960690541682
$ws.Range("D3").Value = 0.04518012690787998
$ws.Range("F3").Value = 6.499322781234298
$ws.Range("G3").Value = 0.002687758981656166
$ws.Range("J3").Value = 0.1758710626823152
$ws.Range("L3").Value = 0.2870734117642968
$ws.Range("M3").Value = 0.5784692023009512

$ws.Range("B4").Value = 3.022285137614347
$ws.Range("C4").Value = 0.1360033035794288
$ws.Range("D4").Value = 0.04121507731058216
$ws.Range("F4").Value = 6.394728955178095
$ws.Range("G4").Value = 0.002692143744943508
$ws.Range("J4").Value = 0.1748879421335694
$ws.Range("L4").Value = 0.2885389116991703
$ws.Range("M4").Value = 0.5769625683503392

$ws.Range("B5").Value = 3.014118118620502
$ws.Range("C5").Value = 0.1327907786742628
$ws.Range("D5").Value = 0.03959703916271451
$ws.Range("F5").Value = 6.35246440323607
$ws.Range("G5").Value = 0.002693984865617241
$ws.Range("J5").Value = 0.1744819989973792
$ws.Range("L5").Value = 0.2891761278176475
$ws.Range("M5").Value = 0.5764606954166993

$ws.Range("B6").Value = 3.01280032173571
$ws.Range("C6").Value = 0.1322600178050379
$ws.Range("D6").Value = 0.03932822146019532
$ws.Range("F6").Value = 6.345467859932768
$ws.Range("G6").Value = 0.002694293867378319
$ws.Range("J6").Value = 0.17441426862597
$ws.Range("L6").Value = 0.289284355565961
$ws.Range("M6").Value = 0.5763841265610949

$ws.Range("B7").Value = 3.022172423866948
$ws.Range("C7").Value = 0.1359597985546657
$ws.Range("D7").Value = 0.04119326534375034
$ws.Range("F7").Value = 6.394157518396725
$ws.Range("G7").Value = 0.002692168354804439
$ws.Range("J7").Value = 0.1748824890883824
$ws.Range("L7").Value = 0.2885473433249359
$ws.Range("M7").Value = 0.5769553462071713

$ws.Range("B8").Value = 3.069110356165652
$ws.Range("C8").Value = 0.1526791573949993
$ws.Range("D8").Value = 0.04940428548670184
$ws.Range("F8").Value = 6.612176876263334
$ws.Range("G8").Value = 0.002683265042539548
$ws.Range("J8").Value = 0.1769022613490883
$ws.Range("L8").Value = 0.2856461882040549
$ws.Range("M8").Value = 0.5804604166775249

$ws.Range("B9").Value = 3.18132002463949
$ws.Range("C9").Value = 0.1869197245910357
$ws.Range("D9").Value = 0.06545664199016699
$ws.Range("F9").Value = 7.051532763981726
$ws.Range("G9").Value = 0.00266750309284811
$ws.Range("J9").Value = 0.1807030835888526
$ws.Range("L9").Value = 0.2812152419151346
$ws.Range("M9").Value = 0.5908820392982435

$ws.Range("B10").Value = 3.276177529124482
$ws.Range("C10").Value = 0.2130054979682257
$ws.Range("D10").Value = 0.07723908065960927
$ws.Range("F10").Value = 7.382002736411209
$ws.Range("G10").Value = 0.002656944475691425
$ws.Range("J10").Value = 0.1834045291089446
$ws.Range("L10").Value = 0.2787238704349093
$ws.Range("M10").Value = 0.6007148955721036

$ws.Range("B11").Value = 3.3220560884821
$ws.Range("C11").Value = 0.2250840118702513
$ws.Range("D11").Value = 0.08260094170053378
$ws.Range("F11").Value = 7.534113163509176
$ws.Range("G11").Value = 0.002652360127845445
$ws.Range("J11").Value = 0.1846151373049203
$ws.Range("L11").Value = 0.27775580827619
$ws.Range("M11").Value = 0.6056637952670272

$ws.Range("B12").Value = 3.339823477083087
$ws.Range("C12").Value = 0.229689058208379
$ws.Range("D12").Value = 0.08463194244488648
$ws.Range("F12").Value = 7.591977015302632
$ws.Range("G12").Value = 0.002650655404259258
$ws.Range("J12").Value = 0.1850710451619726
$ws.Range("L12").Value = 0.2774129484969947
$ws.Range("M12").Value = 0.6076064792610225

$ws.Range("B13").Value = 3.335979388709859
$ws.Range("C13").Value = 0.2286958813206184
$ws.Range("D13").Value = 0.08419449952288005
$ws.Range("F13").Value = 7.579503188574677
$ws.Range("G13").Value = 0.002651021159080447
$ws.Range("J13").Value = 0.1849729677490615
$ws.Range("L13").Value = 0.2774857349834363
$ws.Range("M13").Value = 0.6071850310860896

$ws.Range("B14").Value = 3.323509910796076
$ws.Range("C14").Value = 0.2254622429674384
$ws.Range("D14").Value = 0.08276802001027761
$ws.Range("F14").Value = 7.538868345830565
$ws.Range("G14").Value = 0.002652219253445942
$ws.Range("J14").Value = 0.1846526950692251
$ws.Range("L14").Value = 0.277727125757643
$ws.Range("M14").Value = 0.6058222440848766

$ws.Range("B15").Value = 3.315923388636634
$ws.Range("C15").Value = 0.2234856269469958
$ws.Range("D15").Value = 0.08189434456586753
$ws.Range("F15").Value = 7.514012790595132
$ws.Range("G15").Value = 0.002652957186948066
$ws.Range("J15").Value = 0.1844561934129842
$ws.Range("L15").Value = 0.2778780729779484
$ws.Range("M15").Value = 0.6049964435339774

$ws.Range("B16").Value = 3.27323430683083
$ws.Range("C16").Value = 0.2122204629982889
$ws.Range("D16").Value = 0.07688873539862584
$ws.Range("F16").Value = 7.372098358107451
$ws.Range("G16").Value = 0.002657248459812567
$ws.Range("J16").Value = 0.1833250540666569
$ws.Range("L16").Value = 0.2787904576913647
$ws.Range("M16").Value = 0.6004010671530153

$ws.Range("B17").Value = 3.247745872293024
$ws.Range("C17").Value = 0.2053644285369387
$ws.Range("D17").Value = 0.07381867133480569
$ws.Range("F17").Value = 7.285498966591831
$ws.Range("G17").Value = 0.002659936921924625
$ws.Range("J17").Value = 0.1826265275041088
$ws.Range("L17").Value = 0.2793924765086402
$ws.Range("M17").Value = 0.5977039846553254

$ws.Range("B18").Value = 3.23334219206248
$ws.Range("C18").Value = 0.2014409274522393
$ws.Range("D18").Value = 0.07205301922145679
$ws.Range("F18").Value = 7.235856197979899
$ws.Range("G18").Value = 0.002661503862130727
$ws.Range("J18").Value = 0.1822230252090016
$ws.Range("L18").Value = 0.2797543006760179
$ws.Range("M18").Value = 0.5961974732283011

$ws.Range("B19").Value = 3.228509364613672
$ws.Range("C19").Value = 0.2001158958059079
$ws.Range("D19").Value = 0.0714552195328082
$ws.Range("F19").Value = 7.219076457086516
$ws.Range("G19").Value = 0.002662037946559454
$ws.Range("J19").Value = 0.1820861061770422
$ws.Range("L19").Value = 0.2798794816859598
$ws.Range("M19").Value = 0.5956950789387392

$ws.Range("B20").Value = 3.250432589190893
$ws.Range("C20").Value = 0.2060921996615832
$ws.Range("D20").Value = 0.07414546541539835
$ws.Range("F20").Value = 7.294700292751401
$ws.Range("G20").Value = 0.002659648598907465
$ws.Range("J20").Value = 0.1827010649463716
$ws.Range("L20").Value = 0.2793267806058779
$ws.Range("M20").Value = 0.5979864577325529

$ws.Range("B21").Value = 3.32716178480149
$ws.Range("C21").Value = 0.2264111892304754
$ws.Range("D21").Value = 0.08318699348946268
$ws.Range("F21").Value = 7.550796589748813
$ws.Range("G21").Value = 0.002651866497157882
$ws.Range("J21").Value = 0.184746834537215
$ws.Range("L21").Value = 0.2776555798877425
$ws.Range("M21").Value = 0.6062206629623574

$ws.Range("B22").Value = 3.379607272065812
$ws.Range("C22").Value = 0.2398728708240867
$ws.Range("D22").Value = 0.08909967752525461
$ws.Range("F22").Value = 7.719706357023654
$ws.Range("G22").Value = 0.002646962615240382
$ws.Range("J22").Value = 0.1860692066460885
$ws.Range("L22").Value = 0.2767016193552934
$ws.Range("M22").Value = 0.6120024032412061

$ws.Range("B23").Value = 3.351405139251767
$ws.Range("C23").Value = 0.2326712268520623
$ws.Range("D23").Value = 0.0859435454909061
$ws.Range("F23").Value = 7.629413059352885
$ws.Range("G23").Value = 0.002649563305372027
$ws.Range("J23").Value = 0.1853647371279159
$ws.Range("L23").Value = 0.2771981280267966
$ws.Range("M23").Value = 0.6088798872705397

$ws.Range("B24").Value = 3.249217146059152
$ws.Range("C24").Value = 0.2057631179347652
$ws.Range("D24").Value = 0.07399772380047409
$ws.Range("F24").Value = 7.290539924045902
$ws.Range("G24").Value = 0.002659778883353676
$ws.Range("J24").Value = 0.1826673725264243
$ws.Range("L24").Value = 0.2793564327699585
$ws.Range("M24").Value = 0.5978586143793763

$ws.Range("B25").Value = 3.148793019482696
$ws.Range("C25").Value = 0.1774967109637657
$ws.Range("D25").Value = 0.06111740160096701
$ws.Range("F25").Value = 6.93136175288555
$ws.Range("G25").Value = 0.002671586753848125
$ws.Range("J25").Value = 0.1796913160054956
$ws.Range("L25").Value = 0.2822795474686401
$ws.Range("M25").Value = 0.5876815511680604
